# Adds ci.lower / ci.upper columns (G, H) to the ml_results sheet,
# reflecting a rerun of the analysis scripts that now also reports
# confidence interval bounds alongside the existing statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("G1").Value = "ci.lower"
$ws.Range("H1").Value = "ci.upper"

$ws.Range("G2").Value = 0.158405518943659
$ws.Range("H2").Value = 0.683968552604894

$ws.Range("G3").Value = -0.0491117252880285
$ws.Range("H3").Value = 0.0781731211301816

$ws.Range("G4").Value = -0.0892414740387267
$ws.Range("H4").Value = 0.0493374750514777

$ws.Range("G5").Value = -0.223411654953692
$ws.Range("H5").Value = 0.355613374650606

$ws.Range("G6").Value = -0.37280835438689
$ws.Range("H6").Value = 0.206108461135052

$ws.Range("G7").Value = -0.213900111616524
$ws.Range("H7").Value = 0.340473466103907

$ws.Range("G8").Value = -0.388680323196871
$ws.Range("H8").Value = 0.214883337095084

$ws.Range("G9").Value = -0.0683621159675101
$ws.Range("H9").Value = 0.0629408143949622

$ws.Range("G10").Value = 0.187264191354415
$ws.Range("H10").Value = 0.229464866164573

$ws.Range("G11").Value = 0.0129687030098434
$ws.Range("H11").Value = 0.0559966918195587

# Row 12: no CI was computed (matches blank SE/df/t.ratio/p.value cells already present)
$ws.Cells.Item(12, 7).Style = "Normal"
$ws.Cells.Item(12, 8).Style = "Normal"

# Row 13: no CI was computed (matches blank SE/df/t.ratio/p.value cells already present)
$ws.Cells.Item(13, 7).Style = "Normal"
$ws.Cells.Item(13, 8).Style = "Normal"

# Row 14: no CI was computed (matches blank SE/df/t.ratio/p.value cells already present)
$ws.Cells.Item(14, 7).Style = "Normal"
$ws.Cells.Item(14, 8).Style = "Normal"

$ws.Range("G15").Value = -0.0316463739617981
$ws.Range("H15").Value = 0.0951904672186524

$ws.Range("G16").Value = -0.137831910465116
$ws.Range("H16").Value = 0.414590435247082

$ws.Range("G17").Value = -0.159124364439563
$ws.Range("H17").Value = 0.440776890610853

$ws.Range("G18").Value = -0.0982234505760569
$ws.Range("H18").Value = 0.156346242260363

$ws.Range("G19").Value = -0.427800223233049
$ws.Range("H19").Value = 0.680946932207813

$ws.Range("G20").Value = -0.446823309907384
$ws.Range("H20").Value = 0.711226749301211

$ws.Range("G21").Value = 0.0129687030098434
$ws.Range("H21").Value = 0.0559966918195587

$ws.Range("G22").Value = -0.125881628789924
$ws.Range("H22").Value = 0.13672423193502

$ws.Range("G23").Value = 0.0564835994878654
$ws.Range("H23").Value = 0.243886741101303

$ws.Range("G24").Value = -0.548262034997135
$ws.Range("H24").Value = 0.59548566661154

$ws.Range("G25").Value = 0.0589385086387653
$ws.Range("H25").Value = 0.239963104309987

$ws.Range("G26").Value = -0.554614093318411
$ws.Range("H26").Value = 0.589112266873336
